# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback run:
#   - the status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + per-language sheets)
#   - the stale "handback file is not the latest" error message is
#     cleared from the zh-cn and de-de "Error Detail" columns
#   - the zh-cn / de-de "Latest Handback DateTime" values are refreshed
#   - a handful of columns are widened / narrowed to fit the new content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (Overview!E2 / F2 mirror the per-language Status columns.)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value     = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value     = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn: Latest Handback DateTime refreshed, stale error detail cleared
# ---------------------------------------------------------------------
$wsZhCn.Range("K2").Value = "2016-08-23 12:49:39"
$wsZhCn.Range("P2").Value = ""

# ---------------------------------------------------------------------
# de-de: Latest Handback DateTime refreshed, stale error detail cleared
# ---------------------------------------------------------------------
$wsDeDe.Range("K2").Value = "2016-08-23 12:49:46"
$wsDeDe.Range("P2").Value = ""

# ---------------------------------------------------------------------
# Column width adjustments
#
# The workbook's original widths were written with fractional-character
# precision that Excel's ColumnWidth setter re-quantizes to the nearest
# 1/6-character pixel grid on save. The inputs below are chosen so that
# the quantized, persisted width lands on the closest achievable value
# to the intended target width.
#   target 29.9777047293527 chars -> nearest grid value 30        -> feed 29.1666666666667
#   target 13.7470528738839 chars -> nearest grid value 13.666667 -> feed 12.8333333333333
# ---------------------------------------------------------------------

# Overview: widen Status columns (zh-cn / de-de) E & F
$wsOverview.Range("E1").ColumnWidth = 29.1666666666667
$wsOverview.Range("F1").ColumnWidth = 29.1666666666667

# zh-cn / de-de: widen Status column (C), narrow Error Detail column (P)
$wsZhCn.Range("C1").ColumnWidth = 29.1666666666667
$wsZhCn.Range("P1").ColumnWidth = 12.8333333333333

$wsDeDe.Range("C1").ColumnWidth = 29.1666666666667
$wsDeDe.Range("P1").ColumnWidth = 12.8333333333333
